# Domino Pi BOM: Rev. D -> Rev. E
#  - rename the sheet (and keep Print_Area / print-area history defined
#    names in sync with the new sheet name)
#  - silkscreen font ratio changed to 20% -> the sheet's column widths and
#    two row heights shift by the corresponding small amount
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sheet rename -----------------------------------------------------
$ws.Name = "Domino Pi Rev. E"

# --- column widths (tiny ~0.5% bump from the font-ratio change) -------
$ws.Columns("A:B").ColumnWidth = 4.005882352941177
$ws.Columns("C:C").ColumnWidth = 19.903921568627467
$ws.Columns("D:D").ColumnWidth = 27.249019607843167
$ws.Columns("E:E").ColumnWidth = 30.335294117647067
$ws.Columns("F:F").ColumnWidth = 17.554901960784267
$ws.Columns("G:G").ColumnWidth = 39.32352941176467
$ws.Columns("H:H").ColumnWidth = 61.13921568627446
$ws.Columns("I:I").ColumnWidth = 23.703921568627468

# --- rows 35 & 36 gain an explicit (custom) row height -----------------
$ws.Rows.Item(35).RowHeight = 12.1
$ws.Rows.Item(36).RowHeight = 12.1

# --- print area / named-range history, re-pointed at the new sheet ----
$ws.PageSetup.PrintArea = "A1:I34"
$ws.Names.Add('_xlnm.Print_Area', '=''Domino Pi Rev. E''!$A$2:$I$36')
$ws.Names.Add('_xlnm.Print_Area_0_0_0_0_0_0_0_0_0_0_0_0_0', '=''Domino Pi Rev. E''!$A$1:$I$1')
